# Apply author updates to the "script" sheet of the management/provenance
# comparison workbook, based on replies received from paper authors.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("script")

# Datatrack: Versioning changed from "Trial Sequence" to "Trial Identification"
$ws.Range("D8").Value = "Trial Identification"

# RDataTracker: Versioning changed from "Trial Sequence" to "Trial Identification"
$ws.Range("D18").Value = "Trial Identification"

# RDataTracker: Query updated to reflect interoperable PROV format wording
$ws.Range("G18").Value = "Interoperable Format (PROV), Functions, Proprietary (DDG)"

# SPADE: Query now also mentions SQL/Cypher/Datalog querying in addition to PROV
$ws.Range("G21").Value = "Query (SQL, Cypher, Datalog), Interoperable Format (PROV)"

# versuchung: Query now also mentions SQL querying in addition to Functions
$ws.Range("G26").Value = "Functions, Query (SQL)"
